# Weekly update: a new "Rabanito" price observation for Vega Modelo de
# Temuco is inserted at the top of the data block (row 65), pushing the
# existing rows 65:97 down to 66:98 (the oldest observation, previously
# row 97, survives as the new row 98).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 65, shifting 65:97 -> 66:98.
$ws.Rows("65:65").Insert()

# Populate the new row 65 with this week's observation.
$ws.Cells.Item(65, 1).Value  = 10
$ws.Cells.Item(65, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(65, 3).Value  = "La Araucanía"
$ws.Cells.Item(65, 4).Value  = 45001
$ws.Cells.Item(65, 5).Value  = 9
$ws.Cells.Item(65, 6).Value  = 300000001
$ws.Cells.Item(65, 7).Value  = "Rabanito"
$ws.Cells.Item(65, 8).Value  = "Sin especificar"
$ws.Cells.Item(65, 9).Value  = "Primera"
$ws.Cells.Item(65, 10).Value = 50
$ws.Cells.Item(65, 11).Value = 7000
$ws.Cells.Item(65, 12).Value = 7000
$ws.Cells.Item(65, 13).Value = 7000
$ws.Cells.Item(65, 14).Value = "$/docena de paquetes"
$ws.Cells.Item(65, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(65, 16).Value = 583
$ws.Cells.Item(65, 17).Value = 12
$ws.Cells.Item(65, 18).Value = "Hortaliza"
